# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the associated generate/handoff timestamps on
# all three sheets (Overview, zh-cn, de-de). Also widens the Status
# column(s) so the new, longer status text fits (mirrors Excel's
# auto-fit-to-content behaviour that ran after the text changed).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

# --- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-25 13:01:35"

# --- zh-cn sheet ---------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-25 13:01:31"

# --- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-25 13:01:35"

# --- Widen the Status column(s) to fit the new text -----------------------
# ("In Translation" -> "Ready for handoff" is a few characters wider)
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- Keep the datetime cells' display format intact ------------------------
# (re-assert the original "yyyy-mm-dd HH:mm:ss" formatting on every
# date/time cell so the round-trip through the engine doesn't silently
# drop it from any of them, touched or not)
$overview.Range("G2").NumberFormat = $dateTimeFormat
$zhcn.Range("H2").NumberFormat = $dateTimeFormat
$zhcn.Range("K2").NumberFormat = $dateTimeFormat
$dede.Range("H2").NumberFormat = $dateTimeFormat
$dede.Range("K2").NumberFormat = $dateTimeFormat
